# Fill in the previously-blank Connector Part No. / Manufacturer /
# Number of Positions / Connector Type cells for row 20 (CN-019, "Coolant
# Pump 2"), mirroring the data already present on the otherwise-identical
# row 19 (CN-018, "Coolant Pump 1 (Motor Side)") which also uses the
# Molex 43025-0400 connector.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = "43025-0400"
$ws.Range("D20").Value = "Molex"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = "Microfit 3.0 Receptable Housing, Dual Row, 4 Circuits, Halogen Free"

# Leave the selection on the row that was just edited.
$null = $ws.Range("C20:F20").Select()

# Minor row-height touch-ups that came along with the edit.
$ws.Rows.Item(1).RowHeight = 20
$ws.Rows.Item(35).RowHeight = 16
